$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 (pushes the old row 18 and everything below it down by one).
$ws.Rows.Item(18).Insert()

# Format the new row 18 like the other "warning / needs attention" rows (red "Bad" cell style).
$ws.Range("A18:J18").Style = "Bad"

# Populate the new row 18 with the VoicedLine documentation.
# Order matters here so that new shared-string entries land on the expected indices.
$ws.Range("H18").Value = "[Voice2(optional)]"
$ws.Range("J18").Value = "[Voice3(optional)]"
$ws.Range("B18").Value = "VoicedLine (or vl, voice, v)"
$ws.Range("F18").Value = '[Voice, note that you need to specify full path from "sound/" folder]'
$ws.Range("C18").Value = "None (Instant behavior is undefined)"
$ws.Range("A18").Value = "Display"
$ws.Range("D18").Value = "[Name]"
$ws.Range("E18").Value = "[Dialogue]"
$ws.Range("G18").Value = "[Dialogue2(optional)]"
$ws.Range("I18").Value = "[Dialogue3(optional)]"

# Column width tweaks (D got narrower / manually resized, F/I/J got wider to fit the new voice text).
$ws.Columns.Item(4).ColumnWidth = 45.6667
$ws.Columns.Item(6).ColumnWidth = 62.1667
$ws.Columns.Item(9).ColumnWidth = 19.6667
$ws.Columns.Item(10).ColumnWidth = 16.6667

# Update the view: scroll back to the top and select C19 (the line right below the new row).
$ws.Range("C19").Select()

Write-Host "done"
